# Data retrieved - Thu Jul 29 18:25:29 UTC 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92 had its timestamp recomputed with slightly more precision on re-save
$ws.Cells.Item(92, 1).Value = 44405.76940213078

# Append new row 93 with the latest day's data
$ws.Cells.Item(93, 1).Value = 44406.76769951286
$ws.Cells.Item(93, 1).NumberFormat = $ws.Cells.Item(92, 1).NumberFormat

$ws.Cells.Item(93, 2).Value = 80989
$ws.Cells.Item(93, 3).Value = 68362
$ws.Cells.Item(93, 4).Value = 3687
$ws.Cells.Item(93, 5).Value = 2240
$ws.Cells.Item(93, 6).Value = 1621
$ws.Cells.Item(93, 7).Value = 21319
$ws.Cells.Item(93, 8).Value = 1679
$ws.Cells.Item(93, 9).Value = 904
$ws.Cells.Item(93, 10).Value = 199
